# Linkify email, phone number
#
# 1. "vgrigoriu@gmail.com" -> wrapped in a mailto: hyperlink (Hyperlink style)
# 2. "+4 0751 369 848" -> text corrected to "+40 751 369 848" and wrapped in
#    a tel: hyperlink (Hyperlink style)

$d = $word.ActiveDocument

# --- Email address -> mailto: hyperlink -------------------------------
$emailRange = $d.Content.Duplicate
$emailRange.Find.Execute("vgrigoriu@gmail.com", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Hyperlinks.Add($emailRange, "mailto:vgrigoriu@gmail.com")

# --- Phone number: fix text, then -> tel: hyperlink --------------------
$d.Content.Find.Execute("+4 0751 369 848", $true, $false, $false, $false, $false, $true, 1, $false, "+40 751 369 848", 2)

$phoneRange = $d.Content.Duplicate
$phoneRange.Find.Execute("+40 751 369 848", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Hyperlinks.Add($phoneRange, "tel:+40751369848")

